$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; everything shifts right by one.
$ws.Columns("A:A").Insert()

# New header cell: "Save Label " (new shared string), formatted like the old A1/B1 header cell.
$ws.Range("A1").Value = "Save Label "
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Row 2 sub-header cell in the new column stays blank, matching the old A2 formatting.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Data rows 3-14 get a FALSE flag in the new column, row 15 gets TRUE.
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $false
}
$ws.Cells.Item(15, 1).Value = $true

# The last row's height ticks up slightly once the sheet is relaid out with the extra column.
$ws.Rows.Item(15).RowHeight = 30.75

# Restore the active selection to A16 (below the data), matching the saved state.
$ws.Range("A16").Select() | Out-Null
